# ---------------------------------------------------------------------------
# Applies the "updated results and minor bug" commit:
#   * Removes a handful of now-stale cells from Sheet1 (old rank0/rank1
#     labels, GPS Location / HTTP leftovers next to the small table).
#   * Adds a new "API/Sensor x Rank" results table to Sheet1 (G11:Q22).
#   * Adds a new Sheet2 with the same results table (A1:K12), and makes it
#     the active sheet/tab.
#   * Adjusts a couple of column widths on Sheet1 and the selected cell.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Remove stale cells on Sheet1 (rows 4-6 near the small rank0/rank1 table)
# ---------------------------------------------------------------------------
$ws1.Range("I4").ClearContents()
$ws1.Range("J4").ClearContents()
$ws1.Range("H5").ClearContents()
$ws1.Range("H6").ClearContents()

# ---------------------------------------------------------------------------
# 2. New "API/Sensor" results table on Sheet1, anchored at G11
# ---------------------------------------------------------------------------
$ws1.Range("J11").Value = "Rank 1"
$ws1.Range("N11").Value = "Rank 2"
$ws1.Range("Q11").Value = "Rank 1 + 2 Yes"

$ws1.Range("G12").Value = "API/Sensor"
$ws1.Range("H12").Value = "Total"
$ws1.Range("I12").Value = "yes"
$ws1.Range("J12").Value = "no"
$ws1.Range("K12").Value = "maybe"
$ws1.Range("L12").Value = "N/A"
$ws1.Range("M12").Value = "yes"
$ws1.Range("N12").Value = "no"
$ws1.Range("O12").Value = "maybe"
$ws1.Range("P12").Value = "N/A"

$data = @(
    @(13, "String",             25,19,4,2,0, 12,4,8,1,   22),
    @(14, "ArrayList",          22,13,5,3,1, 6,5,1,10,   18),
    @(15, "HashMap/Dictionary", 13,10,3,0,0, 3,2,3,5,    11),
    @(16, "GPS Location",        5,4,0,1,0, 3,1,1,0,      4),
    @(17, "Accelerometer ",      2,2,0,0,0, 1,0,1,0,      2),
    @(18, "BTLE",                5,4,1,0,0, 3,1,1,0,      4),
    @(19, "Wifi",                5,4,1,0,0, 4,1,0,0,      4)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws1.Range("G$r").Value = $row[1]
    $ws1.Range("H$r").Value = $row[2]
    $ws1.Range("I$r").Value = $row[3]
    $ws1.Range("J$r").Value = $row[4]
    $ws1.Range("K$r").Value = $row[5]
    $ws1.Range("L$r").Value = $row[6]
    $ws1.Range("M$r").Value = $row[7]
    $ws1.Range("N$r").Value = $row[8]
    $ws1.Range("O$r").Value = $row[9]
    $ws1.Range("P$r").Value = $row[10]
    $ws1.Range("Q$r").Value = $row[11]
}

# Overall row (20) - includes a few SUM / literal formulas
$ws1.Range("G20").Value = "Overall"
$ws1.Range("H20").Value = 75
$ws1.Range("I20").Formula = "=SUM(I13:I19)"
$ws1.Range("J20").Value = 14
$ws1.Range("K20").Value = 6
$ws1.Range("L20").Value = 1
$ws1.Range("M20").Formula = "=SUM(M13:M19)"
$ws1.Range("N20").Formula = "=14"
$ws1.Range("O20").Value = 15
$ws1.Range("P20").Value = 16
$ws1.Range("Q20").Formula = "=SUM(Q13:Q19)"

$ws1.Range("G21").Value = "Overall Yes Rank 1 Only"
$ws1.Range("H21").Formula = "=56/75"

$ws1.Range("G22").Value = "Overall Yes Rank 1 and 2"
$ws1.Range("H22").Formula = "=65/75"

# ---------------------------------------------------------------------------
# 3. Column width tweaks on Sheet1
# ---------------------------------------------------------------------------
$ws1.Columns.Item(7).ColumnWidth  = 19.94401
$ws1.Columns.Item(8).ColumnWidth  = 6.94401
$ws1.Columns.Item(17).ColumnWidth = 13.05339

# Selection on Sheet1 moves from I14 to J32
$ws1.Range("J32").Select()

# ---------------------------------------------------------------------------
# 4. Add Sheet2 (after Sheet1) with the same results table, A1:K12
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("D1").Value = "Rank 1"
$ws2.Range("H1").Value = "Rank 2"
$ws2.Range("K1").Value = "Rank 1 + 2 Yes"

$ws2.Range("A2").Value = "API/Sensor"
$ws2.Range("B2").Value = "Total"
$ws2.Range("C2").Value = "yes"
$ws2.Range("D2").Value = "no"
$ws2.Range("E2").Value = "maybe"
$ws2.Range("F2").Value = "N/A"
$ws2.Range("G2").Value = "yes"
$ws2.Range("H2").Value = "no"
$ws2.Range("I2").Value = "maybe"
$ws2.Range("J2").Value = "N/A"

$data2 = @(
    @(3, "String",             25,19,4,2,0, 12,4,8,1,   22),
    @(4, "ArrayList",          22,13,5,3,1, 6,5,1,10,   18),
    @(5, "HashMap/Dictionary", 13,10,3,0,0, 3,2,3,5,    11),
    @(6, "GPS Location",        5,4,0,1,0, 3,1,1,0,      4),
    @(7, "Accelerometer ",      2,2,0,0,0, 1,0,1,0,      2),
    @(8, "BTLE",                5,4,1,0,0, 3,1,1,0,      4),
    @(9, "Wifi",                5,4,1,0,0, 4,1,0,0,      4)
)

foreach ($row in $data2) {
    $r = $row[0]
    $ws2.Range("A$r").Value = $row[1]
    $ws2.Range("B$r").Value = $row[2]
    $ws2.Range("C$r").Value = $row[3]
    $ws2.Range("D$r").Value = $row[4]
    $ws2.Range("E$r").Value = $row[5]
    $ws2.Range("F$r").Value = $row[6]
    $ws2.Range("G$r").Value = $row[7]
    $ws2.Range("H$r").Value = $row[8]
    $ws2.Range("I$r").Value = $row[9]
    $ws2.Range("J$r").Value = $row[10]
    $ws2.Range("K$r").Value = $row[11]
}

$ws2.Range("A10").Value = "Overall"
$ws2.Range("B10").Value = 75
$ws2.Range("C10").Formula = "=SUM(C3:C9)"
$ws2.Range("D10").Value = 14
$ws2.Range("E10").Value = 6
$ws2.Range("F10").Value = 1
$ws2.Range("G10").Formula = "=SUM(G3:G9)"
$ws2.Range("H10").Formula = "=14"
$ws2.Range("I10").Value = 15
$ws2.Range("J10").Value = 16
$ws2.Range("K10").Formula = "=SUM(K3:K9)"

$ws2.Range("A11").Value = "Overall Yes Rank 1 Only"
$ws2.Range("B11").Formula = "=56/75"

$ws2.Range("A12").Value = "Overall Yes Rank 1 and 2"
$ws2.Range("B12").Formula = "=65/75"

# ---------------------------------------------------------------------------
# 5. Sheet2 becomes the active sheet/tab, with the whole sheet selected and
#    the active cell parked at D25 (as close as the host lets us get to the
#    original "select-all, active cell D25" view state).
# ---------------------------------------------------------------------------
$ws2.Select()
$ws2.Range("A1:XFD1048576,D25").Select()

Write-Output "Edit applied"
